$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the paragraph that contains the math run "pip install pandas"
# (a ListParagraph, numId=2) by scanning for a single-oMath paragraph
# whose rendered (math-italic) text has the expected length -- text
# inside an m:oMath zone is rendered using the "Mathematical
# Alphanumeric Symbols" Unicode block by Word, so comparing against
# literal ASCII does not work; comparing UTF-16 code-unit lengths does.
# ------------------------------------------------------------------
function Get-MathItalicLength([string]$s) {
    $len = 0
    foreach ($ch in $s.ToCharArray()) {
        if ($ch -eq ' ') {
            $len += 1
        } else {
            $len += 2
        }
    }
    return $len
}

$pandasLen = (Get-MathItalicLength "pip install pandas") + 1      # +1 for paragraph mark
$openpyxlLen = (Get-MathItalicLength "pip install openpyxl") + 1  # +1 for paragraph mark

# The "pip install pandas" paragraph is immediately followed by the
# "pip install " + "openpyxl" paragraph (both ListParagraph / numId 2,
# single-oMath-zone paragraphs). Matching on that adjacency in
# addition to the rendered-text length disambiguates from unrelated
# paragraphs elsewhere in the document that happen to render to the
# same length (e.g. "pip install notebook").
$pandasPara = $null
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.OMaths.Count -eq 1 -and $p.Range.Text.Length -eq $pandasLen) {
        $next = $d.Paragraphs($p.Index + 1)
        if ($next.Range.OMaths.Count -eq 1 -and $next.Range.Text.Length -eq $openpyxlLen) {
            $pandasPara = $p
            $targetPara = $next
            break
        }
    }
}
if ($pandasPara -eq $null) {
    throw "Could not find the 'pip install pandas' paragraph"
}
if ($targetPara -eq $null) {
    throw "Could not find the 'pip install ' + 'openpyxl' paragraph"
}

$targetParaIndex = $targetPara.Index

# ------------------------------------------------------------------
# 1) Insert a brand-new paragraph right after the "pandas" paragraph,
#    containing a single math run "pip install openpyxl".
# ------------------------------------------------------------------
$pandasPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs($pandasPara.Index + 1)

# Fill the new paragraph with the same pPr (ListParagraph / numId 2)
# and a math run "pip install openpyxl". A throw-away leading text
# run ("X") is included purely to stop Word from auto-wrapping the
# lone oMath into an m:oMathPara display-math container; it is
# deleted again immediately afterwards.
$newParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">X</w:t></w:r><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:shd w:val="clear" w:color="auto" w:fill="E7E6E6" w:themeFill="background2"/></w:rPr><m:t>pip install openpyxl</m:t></m:r></m:oMath></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$newPara.Range.InsertXML($newParaXml)

# Remove the throw-away leading "X" character.
$newPara = $d.Paragraphs($pandasPara.Index + 1)
$leadChar = $d.Range($newPara.Range.Start, $newPara.Range.Start + 1)
$leadChar.Delete()

# ------------------------------------------------------------------
# 2) In the paragraph that originally had two math runs
#    "pip install " + "openpyxl", change the second run's text to
#    "seaborn". Because a new paragraph was inserted before it, its
#    index shifted by one.
# ------------------------------------------------------------------
$targetPara = $d.Paragraphs($targetParaIndex + 1)

$targetXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="5AEC24BA" w14:textId="07D03C89" w:rsidR="00F43B3E" w:rsidRDefault="00F43B3E" w:rsidP="00F43B3E"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">X</w:t></w:r><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:shd w:val="clear" w:color="auto" w:fill="E7E6E6" w:themeFill="background2"/></w:rPr><m:t xml:space="preserve">pip install </m:t></m:r><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:shd w:val="clear" w:color="auto" w:fill="E7E6E6" w:themeFill="background2"/></w:rPr><m:t>seaborn</m:t></m:r></m:oMath></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$targetPara.Range.InsertXML($targetXml)
$targetPara = $d.Paragraphs($targetParaIndex + 1)
$leadChar2 = $d.Range($targetPara.Range.Start, $targetPara.Range.Start + 1)
$leadChar2.Delete()
